# "Generate Report for Handoff"
#
# A fresh handoff was generated for 0cf71a02-1fa4-4dc8-90b5-5c52da7e05a9.md:
#   - its Status flips from "In Translation" to "Ready for handoff"
#   - its "Latest Handoff Datetime" gets a newer timestamp (one per locale)
# Because the report sorts by (status, handoff time) 93237bcb-... (still
# "In Translation") now shows up before 0cf71a02-... on every sheet, i.e.
# rows 2 and 3 swap places. f1a159b1-... (row 4) is untouched.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "93237bcb-bb20-4cdf-8d92-8e150c11553e.md"
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "In Translation"

$ws.Range("A3").Value = "0cf71a02-1fa4-4dc8-90b5-5c52da7e05a9.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

$ws.Range("A4").Value = "f1a159b1-ea13-49ac-adcf-64494b4f3438.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c1bbfe71ec88ad6ab4c2ca421e53862c085cec74/e2e/93237bcb-bb20-4cdf-8d92-8e150c11553e.md", "", "", "93237bcb-bb20-4cdf-8d92-8e150c11553e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2db61dd76e4271feddcb66eb3d5f3a0f59ecf5e7/e2e/0cf71a02-1fa4-4dc8-90b5-5c52da7e05a9.md", "", "", "0cf71a02-1fa4-4dc8-90b5-5c52da7e05a9.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/84646d918171d50de0da9440f0727150ed93a377/e2e/f1a159b1-ea13-49ac-adcf-64494b4f3438.md", "", "", "f1a159b1-ea13-49ac-adcf-64494b4f3438.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/2db61dd76e4271feddcb66eb3d5f3a0f59ecf5e7/.localization-config", "", "", ".localization-config") | Out-Null

# ---- zh-cn sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "93237bcb-bb20-4cdf-8d92-8e150c11553e.md"
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "93237bcb-bb20-4cdf-8d92-8e150c11553e.3f7199c4feff4de859445fd9a385c2e1e000aac5.zh-cn.xlf"
$ws.Range("D2").Value = "2016-02-22 06:09:36"

$ws.Range("A3").Value = "0cf71a02-1fa4-4dc8-90b5-5c52da7e05a9.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "0cf71a02-1fa4-4dc8-90b5-5c52da7e05a9.39666fb79611f01704d27286b6cd994eb2bde93d.zh-cn.xlf"
$ws.Range("D3").Value = "2016-02-22 06:14:12"

$ws.Range("A4").Value = "f1a159b1-ea13-49ac-adcf-64494b4f3438.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "f1a159b1-ea13-49ac-adcf-64494b4f3438.7dd534e680c5590f105c24bb71e86c8fb5b88073.zh-cn.xlf"
$ws.Range("D4").Value = "2016-02-22 06:10:41"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c1bbfe71ec88ad6ab4c2ca421e53862c085cec74/e2e/93237bcb-bb20-4cdf-8d92-8e150c11553e.md", "", "", "93237bcb-bb20-4cdf-8d92-8e150c11553e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7900b988f4ce07bf755b63d07b73a142774d252e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/93237bcb-bb20-4cdf-8d92-8e150c11553e.3f7199c4feff4de859445fd9a385c2e1e000aac5.zh-cn.xlf", "", "", "93237bcb-bb20-4cdf-8d92-8e150c11553e.3f7199c4feff4de859445fd9a385c2e1e000aac5.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2db61dd76e4271feddcb66eb3d5f3a0f59ecf5e7/e2e/0cf71a02-1fa4-4dc8-90b5-5c52da7e05a9.md", "", "", "0cf71a02-1fa4-4dc8-90b5-5c52da7e05a9.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/812541e4ec7db143256de03b2deae23dde58c14c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/0cf71a02-1fa4-4dc8-90b5-5c52da7e05a9.39666fb79611f01704d27286b6cd994eb2bde93d.zh-cn.xlf", "", "", "0cf71a02-1fa4-4dc8-90b5-5c52da7e05a9.39666fb79611f01704d27286b6cd994eb2bde93d.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/84646d918171d50de0da9440f0727150ed93a377/e2e/f1a159b1-ea13-49ac-adcf-64494b4f3438.md", "", "", "f1a159b1-ea13-49ac-adcf-64494b4f3438.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/37edb377889dbb61c151161251a6498350493912/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f1a159b1-ea13-49ac-adcf-64494b4f3438.7dd534e680c5590f105c24bb71e86c8fb5b88073.zh-cn.xlf", "", "", "f1a159b1-ea13-49ac-adcf-64494b4f3438.7dd534e680c5590f105c24bb71e86c8fb5b88073.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/2db61dd76e4271feddcb66eb3d5f3a0f59ecf5e7/.localization-config", "", "", ".localization-config") | Out-Null

# ---- de-de sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "93237bcb-bb20-4cdf-8d92-8e150c11553e.md"
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "93237bcb-bb20-4cdf-8d92-8e150c11553e.3f7199c4feff4de859445fd9a385c2e1e000aac5.de-de.xlf"
$ws.Range("D2").Value = "2016-02-22 06:09:50"

$ws.Range("A3").Value = "0cf71a02-1fa4-4dc8-90b5-5c52da7e05a9.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "0cf71a02-1fa4-4dc8-90b5-5c52da7e05a9.39666fb79611f01704d27286b6cd994eb2bde93d.de-de.xlf"
$ws.Range("D3").Value = "2016-02-22 06:14:25"

$ws.Range("A4").Value = "f1a159b1-ea13-49ac-adcf-64494b4f3438.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "f1a159b1-ea13-49ac-adcf-64494b4f3438.7dd534e680c5590f105c24bb71e86c8fb5b88073.de-de.xlf"
$ws.Range("D4").Value = "2016-02-22 06:10:55"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/c1bbfe71ec88ad6ab4c2ca421e53862c085cec74/e2e/93237bcb-bb20-4cdf-8d92-8e150c11553e.md", "", "", "93237bcb-bb20-4cdf-8d92-8e150c11553e.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/39c11d43d4b3bd4b95754bb725dbfad75fbe1c24/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/93237bcb-bb20-4cdf-8d92-8e150c11553e.3f7199c4feff4de859445fd9a385c2e1e000aac5.de-de.xlf", "", "", "93237bcb-bb20-4cdf-8d92-8e150c11553e.3f7199c4feff4de859445fd9a385c2e1e000aac5.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/2db61dd76e4271feddcb66eb3d5f3a0f59ecf5e7/e2e/0cf71a02-1fa4-4dc8-90b5-5c52da7e05a9.md", "", "", "0cf71a02-1fa4-4dc8-90b5-5c52da7e05a9.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/76ec99e9954a1d9bd631984a06fb73b78713a595/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/0cf71a02-1fa4-4dc8-90b5-5c52da7e05a9.39666fb79611f01704d27286b6cd994eb2bde93d.de-de.xlf", "", "", "0cf71a02-1fa4-4dc8-90b5-5c52da7e05a9.39666fb79611f01704d27286b6cd994eb2bde93d.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/84646d918171d50de0da9440f0727150ed93a377/e2e/f1a159b1-ea13-49ac-adcf-64494b4f3438.md", "", "", "f1a159b1-ea13-49ac-adcf-64494b4f3438.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c7186266aa25d135bb970b3e48207bda8149596e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f1a159b1-ea13-49ac-adcf-64494b4f3438.7dd534e680c5590f105c24bb71e86c8fb5b88073.de-de.xlf", "", "", "f1a159b1-ea13-49ac-adcf-64494b4f3438.7dd534e680c5590f105c24bb71e86c8fb5b88073.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/2db61dd76e4271feddcb66eb3d5f3a0f59ecf5e7/.localization-config", "", "", ".localization-config") | Out-Null

$ws = $wb.Worksheets.Item("Overview")
$ws.Activate()
$ws.Range("A1").Select()
